# Updated CVDs for the month
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Kongegårdsgatan Molndal Sweden" ---
# Row 5 (Internal Fill Rate / Commit-Forecast): clear Jul (O5) value, keep formatting.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("O5").ClearContents()

# --- Sheet 4: "Charlotte  North Carolina" ---
# Row 2 (Professional Voluntary Turnover / Commit-Forecast): update YTD and monthly CVD values.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("E2").Value = 0.7143
$ws4.Range("O2").Value = 0
$ws4.Range("P2").Value = 0.0833333333333333
$ws4.Range("Q2").Value = 0.0833333333333333
$ws4.Range("R2").Value = 0.25
$ws4.Range("S2").Value = 0.0833333333333333
$ws4.Range("T2").Value = 0.0833333333333333
$ws4.Range("U2").Value = 0.0833333333333333
$ws4.Range("V2").Value = 0.25
$ws4.Range("W2").Value = 1

# --- Sheet 9: "Shanghai Minhang District China" ---
# Row 3 (Internal Fill Rate / Commit-Forecast): clear Jul (O3) value, keep formatting.
$ws9 = $wb.Worksheets.Item(9)
$ws9.Range("O3").ClearContents()
